# Adds input/output example data to the workbook:
#  - Лист1 ("мука" example) and Лист2 ("помидор" example) each get extra
#    solver-config columns on row 1 (Группа / Время доставки / В чем дать
#    ответ / other), and their product name (B1) is renamed.
#  - A brand-new Лист3 sheet is added summarising the output (max cost,
#    model choice, optimisation target, priority, grouped products).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Лист1: rename product (meat -> мука) and add new config cells on row 1 ---
$ws1.Activate()
$ws1.Range("B1").Value = "мука"
$ws1.Range("D1").Value = "Группа"
$ws1.Range("E1").Value = 1
$ws1.Range("G1").Value = "Время доставки"
$ws1.Range("H1").Value = 1440
$ws1.Range("J1").Value = "В чем дать ответ"
$ws1.Range("K1").Value = "other"

# --- Лист2: mirror the same new row 1 cells (product name comes later) ---
$ws2.Activate()
$ws2.Range("D1").Value = "Группа"
$ws2.Range("E1").Value = 1
$ws2.Range("G1").Value = "Время доставки"
$ws2.Range("H1").Value = 1440
$ws2.Range("J1").Value = "В чем дать ответ"
$ws2.Range("K1").Value = "other"

# --- Лист3: new sheet with summarised solver output ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Лист3"

$ws3.Range("A1").Value = "Макс. Стоимость"
$ws3.Range("B1").Value = "nan"

$ws3.Range("B2").Value = "Какой моделью решать"
$ws3.Range("C2").Value = "По какому параметру оптимизировать"
$ws3.Range("D2").Value = "Приоритет"
$ws3.Range("E2").Value = "Товары, по которым решение должно проводиться вместе"

$ws3.Range("B3").Value = "fractional"
$ws3.Range("C3").Value = "cost"
$ws3.Range("D3").Value = 1

$ws3.Range("B2:E3").Borders.LineStyle = 1

# Rename Лист2's product last (potato -> помидор) - this is what the
# summary row (E3) below references.
$ws2.Range("B1").Value = "помидор"

$ws3.Range("E3").Value = "мука, помидор"

# Leave behind the final cursor position on each sheet, then come back to
# Лист3 last so it ends up the active tab (matches xr:revisionPtr's
# activeTab="2").
$ws1.Range("A1:K1").Select()
$ws2.Range("H4").Select()

$ws3.Activate()
$excel.ActiveWindow.Zoom = 175
$ws3.Range("F5").Select()
